$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.059167101389184996
$ws.Range("B1").Value = 0.059167100827169181

$ws.Range("A2").Value = 0.074388803205565776
$ws.Range("B2").Value = -0.074388803750090618

$ws.Range("A3").Value = 0.053185770616089971
$ws.Range("B3").Value = -0.053185771188473556

$ws.Range("A4").Value = 0.045777373761806232
$ws.Range("B4").Value = -0.04577737436611453

$ws.Range("A5").Value = -0.032459640022052302
$ws.Range("B5").Value = 0.032459639413312445
